$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# Append 7 new data rows (72-78) below the existing table (rows 1-71), mirroring
# the row-67-71 formatting: column A carries the existing short-date style,
# column F carries the ABS(D-E) temp-diff formula, everything else is a plain
# number or text value (which Excel stores as a shared string).
# ---------------------------------------------------------------------------

# Pull column A's existing date style (s="1") down onto the new date cells by
# copy/paste-special of formats only, so we reuse the workbook's existing
# cellXf instead of minting a new one.
$ws.Range("A71").Copy()
$ws.Range("A72:A78").PasteSpecial(-4122)
[void]($excel.CutCopyMode = 0)

$rows = @(
    @{ Row=72; Date=45797; B="Flowering";    C="Large";  D=54; E=66; G=1.06; H=0.5;  I="Yes"; J=2; K="Bright";  L=3; M=0.46; N=43; O=29.9; P=10; Q=0.94; R=9.9; S=32; T=36 }
    @{ Row=73; Date=45797; B="Nonflowering"; C="Medium"; D=54; E=66; G=1.06; H=0.25; I="Yes"; J=3; K="Bright";  L=3; M=0.46; N=43; O=29.9; P=10; Q=0.94; R=9.9; S=32; T=36 }
    @{ Row=74; Date=45797; B="Nonflowering"; C="Small";  D=54; E=66; G=1.06; H=1;    I="Yes"; J=3; K="Bright";  L=3; M=0.46; N=43; O=29.9; P=10; Q=0.94; R=9.9; S=32; T=36 }
    @{ Row=75; Date=45797; B="Nonflowering"; C="Medium"; D=54; E=66; G=1.06; H=2;    I="Yes"; J=3; K="Bright";  L=3; M=0.46; N=43; O=29.9; P=10; Q=0.94; R=9.9; S=32; T=36 }
    @{ Row=76; Date=45797; B="Nonflowering"; C="Medium"; D=54; E=66; G=1.06; H=3;    I="Yes"; J=3; K="Bright";  L=3; M=0.46; N=43; O=29.9; P=10; Q=0.94; R=9.9; S=32; T=36 }
    @{ Row=77; Date=45797; B="Nonflowering"; C="Large";  D=54; E=66; G=1.06; H=3.5;  I="Yes"; J=4; K="Neutral"; L=3; M=0.46; N=43; O=29.9; P=10; Q=0.94; R=9.9; S=32; T=36 }
    @{ Row=78; Date=45797; B="Tree";         C="Medium"; D=54; E=66; G=1.06; H="=25/3"; I="Yes"; J=1; K="Neutral"; L=3; M=0.46; N=43; O=29.9; P=10; Q=0.94; R=9.9; S=32; T=36 }
)

foreach ($r in $rows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value = $r.Date
    $ws.Cells.Item($row, 2).Value = $r.B
    $ws.Cells.Item($row, 3).Value = $r.C
    $ws.Cells.Item($row, 4).Value = $r.D
    $ws.Cells.Item($row, 5).Value = $r.E
    $ws.Cells.Item($row, 6).Formula = "=ABS(D$row-E$row)"
    $ws.Cells.Item($row, 7).Value = $r.G

    if ($r.H -like "=*") {
        $ws.Cells.Item($row, 8).Formula = $r.H
    } else {
        $ws.Cells.Item($row, 8).Value = $r.H
    }

    $ws.Cells.Item($row, 9).Value = $r.I
    $ws.Cells.Item($row, 10).Value = $r.J
    $ws.Cells.Item($row, 11).Value = $r.K
    $ws.Cells.Item($row, 12).Value = $r.L
    $ws.Cells.Item($row, 13).Value = $r.M
    $ws.Cells.Item($row, 14).Value = $r.N
    $ws.Cells.Item($row, 15).Value = $r.O
    $ws.Cells.Item($row, 16).Value = $r.P
    $ws.Cells.Item($row, 17).Value = $r.Q
    $ws.Cells.Item($row, 18).Value = $r.R
    $ws.Cells.Item($row, 19).Value = $r.S
    $ws.Cells.Item($row, 20).Value = $r.T
}

$excel.Calculate()

# Update the view so the sheet looks like it did right after the paste:
# scrolled down near the bottom of the (now larger) table, with the newly
# appended Visibility column (R72:R78) selected.
[void]$ws.Range("R72:R78").Select()
